$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update selection to D4:D26
$ws.Range("D4:D26").Select()

# Change values in D4:D26 from "Y" to "N"
$ws.Range("D4:D26").Value = "N"

Write-Host "Done"
